$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.592000000000001
$ws.Range("A4").Value = -20.866
$ws.Range("B4").Value = 7.043000000000001
$ws.Range("D4").Value = -7.388999999999998
$ws.Range("B5").Value = 6.250999999999999
$ws.Range("A6").Value = -20.955
$ws.Range("B6").Value = 6.733
$ws.Range("A7").Value = -21.118
$ws.Range("A8").Value = -21.038
$ws.Range("B8").Value = 6.578
$ws.Range("D9").Value = -7.631
$ws.Range("D11").Value = -8.626999999999999
$ws.Range("D14").Value = -7.694999999999999
$ws.Range("A16").Value = -20.727
$ws.Range("B16").Value = 6.759
$ws.Range("D18").Value = -8.323000000000002
$ws.Range("A20").Value = -21.86
$ws.Range("A21").Value = -20.921
$ws.Range("B22").Value = 6.895
$ws.Range("D25").Value = -8.626999999999999
